$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark from its current location
#    (end of the "...se sume en el inventario de productos." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the "Gestión de proveedores" heading and insert a new run
#    " (pendiente)" right after "proveedores", using the same character
#    style (Título 2 Car) as the rest of the heading.
$rng = $d.Content
$null = $rng.Find.Execute("Gestión de proveedores", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" (pendiente)")
$rng.Style = "Título 2 Car"

# The insertion point for the new _GoBack bookmark is right after the
# text we just typed.
$pos = $rng.End

# 3. Re-create the _GoBack bookmark at the new edit location. Word always
#    tracks the last edited spot with a zero-length _GoBack bookmark.
#    We insert a temporary marker character as a separate run, wrap a
#    bookmark around it and then remove the marker, leaving a zero-length
#    bookmark exactly where the edit happened without disturbing the
#    formatting/run structure of the text we just inserted.
$tmp = $d.Range($pos, $pos)
$tmp.InsertAfter("X")
$bmRng = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $bmRng)
$delRng = $d.Range($pos, $pos + 1)
$delRng.Delete()
